$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.7465759548040882, 1.0906706641440915, 1.381463821113775, 2.752351486433011)"
$ws.Range("C2").Value = "NIG(1.3993134252167332, 1.095486443316517, 3.7000130003033926, 5.722441951646958)"
$ws.Range("D2").Value = "JSU(-1.075128692543469, 1.0557342414403257, -0.3273218794053411, 2.5941075992520877)"
$ws.Range("E2").Value = "NIG(1.5802768602754775, 1.2514317073898051, 3.074217655208552, 5.340539688552582)"
